# edit.ps1 — applies the "Fix: Adjust print statements" commit:
#   1) Convert each "Figure <SEQ Figure \* ARABIC>" complex field
#      (fldChar begin/instrText/separate/<num>/end) into the equivalent
#      <w:fldSimple w:instr=" SEQ Figure \* ARABIC "> simple field, for
#      all three figure captions in the document.
#   2) Remove the "To start with, I believe that the spam filter ..."
#      paragraph entirely.
#   3) Move the (hidden) "_GoBack" bookmark from the end of the
#      "Going forward I would continue the project ..." paragraph to its
#      start (i.e. in front of its two runs).

$d = $word.ActiveDocument

# --- 1) Complex SEQ Figure fields -> simple fields ------------------------
# Field.Add() inserts a brand-new <w:fldSimple> right at the supplied
# (collapsed) Range and automatically computes the correct SEQ number from
# context; deleting the original complex field immediately afterwards
# (while it is still present as a placeholder during the Add call) avoids
# disturbing the neighboring text runs. Re-fetch fields by position each
# time since indices shift as we go.
$fieldCount = $d.Fields.Count
for ($i = 1; $i -le $fieldCount; $i++) {
    $oldField = $d.Fields.Item($i)
    $instr = $oldField.Code.Text
    $insertAt = $oldField.Code.Start
    $insertRange = $d.Range($insertAt, $insertAt)

    $d.Fields.Add($insertRange, 1, "SEQ Figure \* ARABIC", $false) | Out-Null

    # The original (now shifted one slot later) complex field still exists;
    # delete it so only the new simple field remains.
    $staleField = $d.Fields.Item($i + 1)
    $staleField.Delete()
}

# --- 2) Delete the "To start with..." paragraph ---------------------------
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("To start with")) {
        $para.Range.Delete()
        break
    }
}

# --- 3) Move the _GoBack bookmark to the start of the "Going forward..." --
#        paragraph (it currently sits at the end of that paragraph).
$paraCount2 = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount2; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Going forward I would continue the project")) {
        $target = $para.Range.Start
        $d.Bookmarks.Add("_GoBack", $d.Range($target, $target)) | Out-Null
        break
    }
}
